# 4.b.1.1.xlsx — add two more year columns (2022, 2023) to the table,
# shrink the footnote row font, and set page setup for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the formatting of column K into the new L/M columns (row by row)
# so that no stray new fonts/styles get created — only the one new
# (smaller, 8pt) font used below for the footnote row is genuinely new.

# Row 3 (thick-bottom-border spacer row)
$ws.Range("K3").Copy()
$ws.Range("L3:M3").PasteSpecial(-4122)

# Row 4 (year headers)
$ws.Range("K4").Copy()
$ws.Range("L4:M4").PasteSpecial(-4122)
$ws.Range("L4").Value = 2022
$ws.Range("M4").Value = 2023

# Row 5
$ws.Range("K5").Copy()
$ws.Range("L5:M5").PasteSpecial(-4122)
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = 700

# Row 6 (K6 style carries to both L6 and M6)
$ws.Range("K6").Copy()
$ws.Range("L6:M6").PasteSpecial(-4122)
$ws.Range("L6").Value = "-"
$ws.Range("M6").Value = 6

# Row 7 (L7 matches K7's own style; M7 instead matches the K6/K8 style)
$ws.Range("K7").Copy()
$ws.Range("L7").PasteSpecial(-4122)
$ws.Range("L7").Value = 23

$ws.Range("K6").Copy()
$ws.Range("M7").PasteSpecial(-4122)
$ws.Range("M7").Value = "-"

# Row 8
$ws.Range("K8").Copy()
$ws.Range("L8:M8").PasteSpecial(-4122)
$ws.Range("L8").Value = 7
$ws.Range("M8").Value = 5

# Row 9
$ws.Range("K9").Copy()
$ws.Range("L9:M9").PasteSpecial(-4122)
$ws.Range("L9").Value = 23
$ws.Range("M9").Value = 21

# Row 10
$ws.Range("K10").Copy()
$ws.Range("L10:M10").PasteSpecial(-4122)
$ws.Range("L10").Value = 172
$ws.Range("M10").Value = 143

$excel.CutCopyMode = 0

# --- Footnote row (11): shrink the font from 9pt to 8pt (creates the new font/style)
$ws.Range("A11:C11").Font.Size = 8

# --- Page setup for printing
$ws.PageSetup.PaperSize = 256
$ws.PageSetup.Orientation = 1
